$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four new "Template" related header/value cells (Z2, AA2, Z3, Z4) ---
# Set values first.
$ws.Range("Z2").Value = "Name of the Template"
$ws.Range("AA2").Value = "Templates and its location with full path"
$ws.Range("Z3").Value = "Location of the template in the particular folder"
$ws.Range("Z4").Value = "Pdf Location of the Template (path of the Pdf file)"

# Apply the full style (small font, thin box border, centered+wrapped text) to Z2 first.
$z2 = $ws.Range("Z2")
$z2.Font.Size = 7.5
$z2.Borders.LineStyle = 1
$z2.Borders.Weight = 2
$z2.VerticalAlignment = -4108
$z2.WrapText = $true

# Copy that exact style onto the other three cells (one at a time, so the
# same cell style / xf index gets reused for all of them).
$z2.Copy()
$ws.Range("AA2").PasteSpecial(-4122)
$z2.Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$z2.Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# AA3 stays empty but gets a wrap-text only style.
$ws.Range("AA3").WrapText = $true

# --- Column widths for the two new columns ---
$ws.Columns("Z").ColumnWidth = 43.333333333333336
$ws.Columns("AA").ColumnWidth = 42

# --- Update the visible selection / active cell to match the edited area ---
$ws.Activate()
$ws.Range("Z2:AA4").Select()
